# Apply the changes described in the commit: add a new "VariableCostsBudget"
# worksheet with a budgeted-variable-costs calculation, and extend the
# "RicaviBudget" sheet with a new "Famiglia" (P/V/A) classification column
# (plus a couple of explanatory notes) that the new sheet's SUMIF formulas
# rely on.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new worksheet between "VariableCosts" and "FixedCosts"
# ---------------------------------------------------------------------
$after = $wb.Worksheets.Item("VariableCosts")
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$newWs.Name = "VariableCostsBudget"

$rb = $wb.Worksheets.Item("RicaviBudget")

# ---------------------------------------------------------------------
# 2) Populate cells in the exact order the original author typed them, so
#    that the shared-string table in sharedStrings.xml ends up built in
#    the same order (new unique strings are appended in first-use order).
# ---------------------------------------------------------------------

# 2a. First few labels of the new sheet
$newWs.Range("E3").Value = "Voce di costo"
$newWs.Range("H8").Value = "Formula/base di cacolo"
$newWs.Range("G8").Value = "Valore Budget (€)"

# 2b. Switch to RicaviBudget: add "Famiglia" column header + P/V/A values
$rb.Range("G2").Value = "Famiglia"
$rb.Range("G22:G28").Value = "P"
$rb.Range("G29:G44").Value = "V"
$rb.Range("G45:G63").Value = "A"

# 2c. Explanatory notes for the first row of each family group
$rb.Range("H22").Value = "Polveri P, comsumi energia/metano diversi"
$rb.Range("H29").Value = "Polveri V, comsumi medi diversi"
$rb.Range("H45").Value = "Polveri A, consumi diversi dai primi due"

# 2d. "Note" header moves from G2 to H2 (same string, already existing)
$rb.Range("H2").Value = "Note"

# 2e. Move the two existing note formulas from column G to column H for
#     rows 64 and 65 (column G is now used for the family classification)
$rb.Range("H64").Formula = $rb.Range("G64").Formula
$rb.Range("H65").Formula = $rb.Range("G65").Formula
$rb.Range("G64").ClearContents()
$rb.Range("G65").ClearContents()

Write-Host "stage2 ok"

# ---------------------------------------------------------------------
# 3) Back to the new sheet: finish the E/F/G mini-table (rows 5-7)
# ---------------------------------------------------------------------
$newWs.Range("E5").Value = "PF V"
$newWs.Range("E6").Value = "PF A"
$newWs.Range("E7").Value = "Totale PF"
$newWs.Range("F3").Value = "quantita budget (kg)"
$newWs.Range("G3").Value = "Quantità Budget (ton)"

Write-Host "stage3 ok"

# ---------------------------------------------------------------------
# 4) Labels down column B (rows 2-6, then 8-9, 11-12, 14-15, 17-18)
# ---------------------------------------------------------------------
$newWs.Range("B2").Value = "Costo energia totale 2021"
$newWs.Range("C2").Value = 3450000
$newWs.Range("B3").Value = "% elettricità"
$newWs.Range("C3").Value = 0.41
$newWs.Range("B4").Value = "% metano"
$newWs.Range("C4").Formula = "=1-C3"
$newWs.Range("B5").Value = "kWh 2021"
$newWs.Range("C5").Value = 8997222.22
$newWs.Range("B6").Value = "smc metano 2021"
$newWs.Range("C6").Value = 4300000

Write-Host "stage4 ok"

# ---------------------------------------------------------------------
# 5) Family quantities table (rows 4-7), referencing RicaviBudget via
#    SUMIF on the new "Famiglia" column
# ---------------------------------------------------------------------
$newWs.Range("E4").Value = "PF P"
$newWs.Range("F4").Formula = "=SUMIF(RicaviBudget!`$G:`$G,""P"",RicaviBudget!`$D:`$D)"
$newWs.Range("G4").Formula = "=F4/1000"

$newWs.Range("F5").Formula = "=SUMIF(RicaviBudget!`$G:`$G,""V"",RicaviBudget!`$D:`$D)"
$newWs.Range("F6").Formula = "=SUMIF(RicaviBudget!`$G:`$G,""A"",RicaviBudget!`$D:`$CD)"
$newWs.Range("G5:G6").Formula = "=F5/1000"

$newWs.Range("F7").Formula = "=SUM(F4:F6)"
$newWs.Range("G7").Formula = "=G6+G5+G4"

Write-Host "stage5 ok"

# ---------------------------------------------------------------------
# 6) Row 8: electricity cost split + the F8/G8/H8 mini-table header
#    (G8/H8 were already set in stage 1)
# ---------------------------------------------------------------------
$newWs.Range("B8").Value = "Costo elettricità 2021"
$newWs.Range("C8").Formula = "=C2*C3"
$newWs.Range("F8").Value = "Note"

$newWs.Range("B9").Value = "Costo Metano 2021"
$newWs.Range("C9").Formula = "=C2*C4"

Write-Host "stage6 ok"

# ---------------------------------------------------------------------
# 7) Price-per-unit rows (11-12) and increase-rate rows (14-15)
# ---------------------------------------------------------------------
$newWs.Range("B11").Value = "Prezzo €/kWh 2021"
$newWs.Range("C11").Formula = "=C8/C5"

$newWs.Range("B12").Value = "Prezzo €/smc 2021"
$newWs.Range("C12").Formula = "=C9/C6"

$newWs.Range("B14").Value = "Aumento Elettricità 2022"
$newWs.Range("C14").Value = 0.11

$newWs.Range("B15").Value = "Aumento Metano 2022"
$newWs.Range("C15").Value = 0.14

$newWs.Range("B17").Value = "Prezzo €/kWh 2022"
$newWs.Range("C17").Formula = "=C11*(1+C14)"

$newWs.Range("B18").Value = "Prezzo €/smc 2022"
$newWs.Range("C18").Formula = "=C12*(1+C15)"

Write-Host "stage7 ok"

# ---------------------------------------------------------------------
# 8) Final cost-summary table (rows 9-15), filled column-major within
#    each row: H notes for rows 9-10 first, then the F labels 9-15
# ---------------------------------------------------------------------
$newWs.Range("G9").Formula = "=(G6*0.93+G5*0.64)*C18"
$newWs.Range("H9").Value = "=(ton PF A*0,93 + ton PF V*0,64)*prezzo smc 2022"

$newWs.Range("G10").Formula = "=(G4*140 + G5*60 + G6*65 + 3300000) * C17"
$newWs.Range("H10").Value = "=(ton P*140 + ton V*60 + ton A*65 + 3.300.000)*prezzo kWh 2022"

$newWs.Range("F9").Value = "Costo metano"
$newWs.Range("F10").Value = "Costo energia elettrica"

$newWs.Range("G11").Formula = "=-'CE 21-22'!F15*1.025"
$newWs.Range("F11").Value = "Materiali di consumo"

$newWs.Range("G12").Formula = "=G7*10.24"
$newWs.Range("F12").Value = "Pulizia e smaltimento rifiuti"

$newWs.Range("G13").Formula = "=G7*43.2"
$newWs.Range("F13").Value = "Trasporti vendita"

$newWs.Range("G14").Formula = "=(SUMIF(RicaviBudget!`$G:`$G,""P"",RicaviBudget!`$F:`$F)+SUMIF(RicaviBudget!`$G:`$G,""V"",RicaviBudget!`$F:`$F))*0.02"
$newWs.Range("F14").Value = "Provvigioni su vendite"

$newWs.Range("G15").Formula = "=SUM(G9:G14)"
$newWs.Range("F15").Value = "Totale costi variabili"

Write-Host "stage8 ok"
